$d = $word.ActiveDocument

# Helper: replace a Range's contents with a raw OOXML fragment of one or
# more <w:p> block-level elements (InsertXML REPLACES the range contents).
function Set-RangeXml($range, [string]$innerXml) {
    $wrapper = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData>' +
               '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:body>' + $innerXml + '</w:body>' +
               '</w:document>' +
               '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($wrapper)
}

# --- 1) "Drop missing values" paragraph: add <w:ilvl w:val="0"/> and split
#        the run that carried both the <w:tab/> and the text into two runs.
$p17 = $d.Paragraphs(17)
$xml17 = '<w:p>' +
           '<w:pPr>' +
             '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr>' +
             '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
           '</w:pPr>' +
           '<w:r>' +
             '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
             '<w:tab/>' +
           '</w:r>' +
           '<w:r>' +
             '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
             '<w:t xml:space="preserve">Drop missing values </w:t>' +
           '</w:r>' +
         '</w:p>'
Set-RangeXml $p17.Range $xml17

# --- 2) "For columns -> df.drop()" paragraph: add <w:ilvl w:val="0"/>.
$p18 = $d.Paragraphs(18)
$xml18 = '<w:p>' +
           '<w:pPr>' +
             '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr>' +
             '<w:ind w:firstLine="720" w:firstLineChars="0"/>' +
             '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
           '</w:pPr>' +
           '<w:r>' +
             '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
             '<w:t>For columns -&gt; df.drop()</w:t>' +
           '</w:r>' +
         '</w:p>'
Set-RangeXml $p18.Range $xml18

# --- 3) "For rows df.dropna()" paragraph: add <w:ilvl w:val="0"/>.
$p19 = $d.Paragraphs(19)
$xml19 = '<w:p>' +
           '<w:pPr>' +
             '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr>' +
             '<w:ind w:firstLine="720" w:firstLineChars="0"/>' +
             '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
           '</w:pPr>' +
           '<w:r>' +
             '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
             '<w:t xml:space="preserve">For rows df.dropna() </w:t>' +
           '</w:r>' +
         '</w:p>'
Set-RangeXml $p19.Range $xml19

# --- 4) "Dealing with null and empty values" paragraph: add <w:ilvl w:val="0"/>.
$p20 = $d.Paragraphs(20)
$xml20 = '<w:p>' +
           '<w:pPr>' +
             '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr>' +
             '<w:ind w:firstLine="720" w:firstLineChars="0"/>' +
             '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
           '</w:pPr>' +
           '<w:r>' +
             '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
             '<w:t>Dealing with null and empty values</w:t>' +
           '</w:r>' +
         '</w:p>'
Set-RangeXml $p20.Range $xml20

# --- 5) "Imputataion" paragraph: add <w:ilvl w:val="0"/>.
$p21 = $d.Paragraphs(21)
$xml21 = '<w:p>' +
           '<w:pPr>' +
             '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr>' +
             '<w:ind w:firstLine="720" w:firstLineChars="0"/>' +
             '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
           '</w:pPr>' +
           '<w:r>' +
             '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
             '<w:t>Imputataion</w:t>' +
           '</w:r>' +
         '</w:p>'
Set-RangeXml $p21.Range $xml21

# --- 6) "Advanced imputation: " paragraph: add <w:ilvl w:val="0"/>.
$p24 = $d.Paragraphs(24)
$xml24 = '<w:p>' +
           '<w:pPr>' +
             '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr>' +
             '<w:ind w:firstLine="720" w:firstLineChars="0"/>' +
             '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
           '</w:pPr>' +
           '<w:r>' +
             '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
             '<w:t xml:space="preserve">Advanced imputation: </w:t>' +
           '</w:r>' +
         '</w:p>'
Set-RangeXml $p24.Range $xml24

# --- 7) "SMOTE" paragraph loses the _GoBack bookmark (it moves further
#        down), and gets four new paragraphs inserted right after it: a
#        page break, a "Feature Engineering and  Selection" heading, a
#        "Process of creating features..." paragraph (carrying the
#        relocated bookmark), and two blank paragraphs.
$p26 = $d.Paragraphs(26)
$smotePara = '<w:p>' +
               '<w:pPr>' +
                 '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>' +
                 '<w:ind w:firstLine="720" w:firstLineChars="0"/>' +
                 '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
               '</w:pPr>' +
               '<w:r>' +
                 '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
                 '<w:t>SMOTE</w:t>' +
               '</w:r>' +
               '<w:r>' +
                 '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
                 '<w:tab/>' +
               '</w:r>' +
             '</w:p>'
$pageBreakPara = '<w:p>' +
                    '<w:pPr>' +
                      '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
                    '</w:pPr>' +
                    '<w:r>' +
                      '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
                      '<w:br w:type="page"/>' +
                    '</w:r>' +
                  '</w:p>'
$headingPara = '<w:p>' +
                 '<w:pPr>' +
                   '<w:rPr>' +
                     '<w:rFonts w:hint="default"/><w:b/><w:bCs/>' +
                     '<w:sz w:val="32"/><w:szCs w:val="32"/>' +
                     '<w:lang w:val="en-US"/>' +
                   '</w:rPr>' +
                 '</w:pPr>' +
                 '<w:r>' +
                   '<w:rPr>' +
                     '<w:rFonts w:hint="default"/><w:b/><w:bCs/>' +
                     '<w:sz w:val="32"/><w:szCs w:val="32"/>' +
                     '<w:lang w:val="en-US"/>' +
                   '</w:rPr>' +
                   '<w:t>Feature Engineering and  Selection</w:t>' +
                 '</w:r>' +
               '</w:p>'
$processPara = '<w:p>' +
                 '<w:pPr>' +
                   '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
                 '</w:pPr>' +
                 '<w:r>' +
                   '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
                   '<w:t>Process of creating features that enhance the performance of ML models</w:t>' +
                 '</w:r>' +
                 '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
                 '<w:bookmarkEnd w:id="0"/>' +
               '</w:p>'
$emptyPara = '<w:p>' +
               '<w:pPr>' +
                 '<w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr>' +
               '</w:pPr>' +
             '</w:p>'

$xml26 = $smotePara + $pageBreakPara + $headingPara + $processPara + $emptyPara + $emptyPara
Set-RangeXml $p26.Range $xml26
